$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "cours 8 - 9" -> mark "Laboratoire 6" (column J) as "Ok" for every
# student row (rows 2-16 and 18-27; row 17/28 stay untouched, rows
# 29-31 are summary/formula rows and are not touched here).
$rows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,18,19,20,21,22,23,24,25,26,27)
foreach ($r in $rows) {
    $ws.Range("J$r").Value = "Ok"
}

# Freeze the first two columns (Prénom / Nom) and scroll the view so the
# newly-filled "Laboratoire 6" column is visible, leaving J28 selected.
$ws.Activate()
$aw = $excel.ActiveWindow
$ws.Range("C1").Select()
$aw.FreezePanes = $true
$ws.Range("J28").Select()
